# Auto-generated edit script: rebuild PMI sections with full_text tweets,
# limited to 5 most recent tweets per country, reordered as
# China, US, Australia, UK, Germany, Japan.

$d = $word.ActiveDocument

# Base64-encoded "Style|Text" entries (UTF-8) -- avoids any PowerShell
# quoting/escaping issues with the tweet text (quotes, $, emoji, etc.)
$items = @(
    "Heading2|LS0tLS0tLS0tLS0tLSBDaGluYSBQTUkgLS0tLS0tLS0tLS0tLS0tLS0tLS0tIDo=",
    "Normal|MjAyMy0wNC0wNiAwODoyOTowNyswMDowMPCfh6jwn4ezICNDaGluYSdzIHNlcnZpY2Ugc2VjdG9yIHJlZ2lzdGVyZWQgaXRzIHN0cm9uZ2VzdCBleHBhbnNpb24gc2luY2UgTm92IDIwMjAgaW4gTWFyY2ggKCNQTUk6IDU3Ljg7IEZlYjogNTUuMCksIGhlbHBlZCBieSBhIHJlY29yZCBpbmNyZWFzZSBpbiBuZXcgI2V4cG9ydCBidXNpbmVzcy4gQ29zdCBwcmVzc3VyZXMgcGlja2VkIHVwIHRvIHRoZSBzdGVlcGVzdCBmb3IgNyBtb250aHMuIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL0t2eGxxaGl3ZEQgaHR0cHM6Ly90LmNvL3VnMVg0VHoyR1U=",
    "Normal|MjAyMy0wNC0wMyAwMjowNjo0MCswMDowMFRoZSBDYWl4aW4gQ2hpbmEgR2VuZXJhbCBNYW51ZmFjdHVyaW5nIFBNSSBwb3N0ZWQgNTAuMCBpbiBNYXJjaCAoRmViOiA1MS42KS4gTWFudWZhY3R1cmluZyBwcm9kdWN0aW9uIGFuZCBuZXcgb3JkZXJzIHJlZ2lzdGVyZWQgc29mdGVyIGdyb3d0aCB3aGlsZSBzdXBwbHkgY2hhaW5zIGltcHJvdmVkIGZvciBhIHNlY29uZCBtb250aCBydW5uaW5nLiBAY2FpeGluIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL3dVZ1BWWGtXbG4gaHR0cHM6Ly90LmNvL3NYSmlyTjlmbzI=",
    "Normal|MjAyMy0wMy0wMyAwMTo0OToxNCswMDowMENoaW5hJ3MgYnVzaW5lc3MgYWN0aXZpdHkgcm9zZSBzaGFycGx5IGluIEZlYnJ1YXJ5IGFjY29yZGluZyB0byB0aGUgQ2FpeGluIENoaW5hIEdlbmVyYWwgU2VydmljZXMgI1BNSSwgd2hpY2ggcm9zZSB0byA1NS4wIChKYW46IDUyLjkpLiBTZXJ2aWNlcyBjb21wYW5pZXMgc2VlIHN0ZWVwZXIgaW5jcmVhc2VzIGluIGFjdGl2aXR5IGFuZCBuZXcgd29yayBhbWlkIG11dGVkIGluZmxhdGlvbmFyeSBwcmVzc3VyZXMuIEBjYWl4aW4gUmVhZCBtb3JlOiBodHRwczovL3QuY28vNk9MZHp1V1dxNyBodHRwczovL3QuY28vOVhwcG9nZVdPaA==",
    "Normal|MjAyMy0wMy0wMSAwMTo0ODo0MSswMDowMENoaW5hJ3MgbWFudWZhY3R1cmluZyBzZWN0b3IgcmV0dXJuZWQgdG8gZ3Jvd3RoIGluIEZlYnJ1YXJ5IHdpdGggdGhlIENhaXhpbiBDaGluYSBHZW5lcmFsIE1hbnVmYWN0dXJpbmcgI1BNSSBhdCA1MS42IChKYW46IDQ5LjIpLiBTdXBwbGllcnMnIGRlbGl2ZXJ5IHRpbWVzIGltcHJvdmVkIGF0IHRoZSBxdWlja2VzdCByYXRlIGZvciBlaWdodCB5ZWFycy4gQGNhaXhpbiAgUmVhZCBtb3JlOiBodHRwczovL3QuY28vWFlBdlJ3OHhTSSBodHRwczovL3QuY28vZTNVS3hnWUo1Ug==",
    "Normal|MjAyMy0wMi0xMyAwOTo0ODoyOCswMDowMFR1bmUgaW50byBvdXIgZmlyc3QgUE1JLWJhc2VkIHBvZGNhc3Qgb2YgMjAyMyB0byBoZWFyIGFib3V0IHdoYXQgaXMgaGFwcGVuaW5nIGluIGdsb2JhbCAjc3VwcGx5IGNoYWlucywgcGx1cyBhbiB1cGRhdGUgb24gI2luZmxhdGlvbiBpbiB0aGUgZXVyb3pvbmUgYW5kIHRoZSBpbXBhY3Qgb2YgbG9vc2VyICNwYW5kZW1pYyByZXN0cmljdGlvbnMgaW4gTWFpbmxhbmQgQ2hpbmE6IGh0dHBzOi8vdC5jby9VbzhQVWxCeWpIICgxLzkp",
    "Heading2|LS0tLS0tLS0tLS0tLSBVUyBQTUkgLS0tLS0tLS0tLS0tLS0tLS0tLS0tIDo=",
    "Normal|MjAyMy0wNC0wMyAxMzo1MzozOCswMDowMPCfh7rwn4e4IFRoZSAjVVMgbWFudWZhY3R1cmluZyBzZWN0b3Igc2lnbmFsbGVkIGEgbWlsZCBkZXRlcmlvcmF0aW9uIGluIG9wZXJhdGluZyBjb25kaXRpb25zIGR1cmluZyBNYXJjaCAoI1BNSSBhdCA0OS4yOyBGZWIgNDcuMykuIFdoaWxlIG91dHB1dCByb3NlIGZvciB0aGUgZmlyc3QgdGltZSBpbiBmaXZlIG1vbnRocywgbXV0ZWQgY2xpZW50IGRlbWFuZCBsZWQgdG8gYSBmdXJ0aGVyIGRlY2xpbmUgaW4gbmV3IG9yZGVycy4gUmVhZCBtb3JlOiBodHRwczovL3QuY28vbHA0OWNBQXlPbyBodHRwczovL3QuY28vNzdBT2RnMmVCMQ==",
    "Normal|MjAyMy0wMy0yNCAxNTowMDowMyswMDowMPCfh7rwn4e4IEF0IDUzLjMgaW4gTWFyY2ggKEZlYjogNTAuMSksIGZsYXNoICNQTUkgZmlndXJlcyBmb3IgdGhlIFVTIHNpZ25hbGxlZCBzb2xpZCBncm93dGggYWNyb3NzIHRoZSBwcml2YXRlIHNlY3Rvci4gSW1wcm92ZWQgZGVsaXZlcnkgdGltZXMgYXQgZ29vZHMgcHJvZHVjZXJzIGFuZCBzZXJ2aWNlIHByb3ZpZGVycyBzZWVpbmcgYW4gdXB0aWNrIGluIGRlbWFuZCBjb25kaXRpb25zIGhlbHBlZCBkcml2ZSB0b3RhbCBhY3Rpdml0eS4gUmVhZCBtb3JlOiBodHRwczovL3QuY28vYUptelhPdUZaMSBodHRwczovL3QuY28vRUhCUXV1T0pLUg==",
    "Normal|MjAyMy0wMy0wMSAxNDo1ODo1OSswMDowMPCfh7rwn4e4IFdoaWxlIHRoZSBkZWNsaW5lIGFjcm9zcyB0aGUgI1VTIG1hbnVmYWN0dXJpbmcgc2VjdG9yIHJlbWFpbmVkIHNvbGlkLCB0aGUgI1BNSSBoaXQgYSB0aHJlZS1tb250aCBoaWdoIG9mIDQ3LjMgaW4gRmVicnVhcnkgKEphbjogNDYuOSksIGluZGljYXRpbmcgYSBzb2Z0ZXIgY29udHJhY3Rpb24uIE1vcmVvdmVyLCBpbXByb3ZlbWVudCBpbiBzdXBwbHkgY2hhaW5zIGhlbHBlZCByZWR1Y2UgaW5wdXQgY29zdCBpbmZsYXRpb24uIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL0Q4RHNraERWdHAgaHR0cHM6Ly90LmNvL1MwT21MR2JydHM=",
    "Normal|MjAyMy0wMi0yMSAxNDo1Nzo0NiswMDowMPCfh7rwn4e4IEZlYnJ1YXJ5IGZsYXNoICNQTUkgcmVhZGluZyBmb3IgdGhlICNVUyBwcml2YXRlIHNlY3RvciBzaWduYWxsZWQgYSBicm9hZCBzdGFibGlzYXRpb24sIGFzIHRoZSBQTUkgcm9zZSB0byBhbiBlaWdodC1tb250aCBoaWdoIG9mIDUwLjIuIEdyb3d0aCB3YXMgY29uZmluZWQgdG8gdGhlIHNlcnZpY2Ugc2VjdG9yLCB3aGlsZSB0aGUgZG93bnR1cm4gaW4gb3V0cHV0IGFjcm9zcyBnb29kcyBwcm9kdWNlcnMgY29udGludWVkLiBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9LUWtEWFR5ZHdtIGh0dHBzOi8vdC5jby9XZnp3Mm5HZ1M4",
    "Normal|MjAyMy0wMi0wMyAxNTozNjowNSswMDowMPCfh7rwn4e4IFRoZSBVUyBzZXJ2aWNlIHNlY3RvciBzaWduYWxsZWQgYSBmdXJ0aGVyIGNvbnRyYWN0aW9uICBkdXJpbmcgSmFudWFyeSBhcyB3ZWFrIGNsaWVudCBjb25kaXRpb25zIGhhbXBlcmVkIGJ1c2luZXNzIGFjdGl2aXR5LiBIb3dldmVyLCB0aGUgbGF0ZXN0ICNQTUkgcmVhZGluZyB0aWNrZWQgdXAgdG8gNDYuOCAoRGVjOiA0NC43KSB0byBzaWduYWwgdGhlIHNvZnRlc3QgcmVkdWN0aW9uIGluIHRocmVlIG1vbnRocy4gUmVhZCBtb3JlOiBodHRwczovL3QuY28vclRwd2JCQWZ4NSBodHRwczovL3QuY28vMDByN3BYZGF4Mw==",
    "Heading2|LS0tLS0tLS0tLS0tLSBBdXN0cmFsaWEgUE1JIC0tLS0tLS0tLS0tLS0tLS0tLS0tLSA6",
    "Normal|MjAyMy0wNC0wNSAwMDowMjo0NyswMDowMEF1c3RyYWxpYeKAmXMgc2VydmljZSBzZWN0b3IgcmVnaXN0ZXJlZCBtaWxkIGNvbnRyYWN0aW9uIGluIE1hcmNoIHdpdGggdGhlIEp1ZG8gQmFuayBBdXN0cmFsaWEgU2VydmljZXMgI1BNSSBwb3N0aW5nIDQ4LjYgKEZlYjogNTAuNykuIEJ1c2luZXNzIGNvbmZpZGVuY2UgaW1wcm92ZWQsIGhvd2V2ZXIuIEBKdWRvX0JhbmsgUmVhZCBtb3JlOiBodHRwczovL3QuY28vcGs3ZVRhaHUzaSBodHRwczovL3QuY28vSmhYMFJFNDNmTQ==",
    "Normal|MjAyMy0wNC0wMyAwMjoxMTozMyswMDowMEF1c3RyYWxpYSdzIG1hbnVmYWN0dXJpbmcgc2VjdG9yIHJlZ2lzdGVyZWQgbWlsZCBjb250cmFjdGlvbiBpbiBNYXJjaCB3aXRoIHRoZSBKdWRvIEJhbmsgQXVzdHJhbGlhIE1uYXVmYWN0dXJpbmcgI1BNSSBwb3N0aW5nIDQ5LjEgKEZlYjogNTAuNSkuIFByaWNlIHByZXNzdXJlcyBlYXNlZCwgaG93ZXZlci4gQEp1ZG9fQmFuayBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9WS2tOMnZpU3Z0IGh0dHBzOi8vdC5jby9nc0NiQXNUbXV3",
    "Normal|MjAyMy0wMy0yNCAwMDowNDozOCswMDowMEF1c3RyYWxpYeKAmXMgcHJpdmF0ZSBzZWN0b3IgYWN0aXZpdHkgc29mdGVuZWQgaW4gTWFyY2ggd2l0aCB0aGUgSnVkbyBCYW5rIEZsYXNoIEF1c3RyYWxpYSBDb21wb3NpdGUgI1BNSSBwb3N0aW5nIDQ4LjEgKEZlYiBmaW5hbDogNTAuNikuIFRoYXQgc2FpZCwgYnVzaW5lc3MgY29uZmlkZW5jZSBpbXByb3ZlZC4gQEp1ZG9fQmFuayBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby8zZ0JoNzBiOU9hIGh0dHBzOi8vdC5jby9DQzBiUVhwS0Yw",
    "Normal|MjAyMy0wMy0wMiAyMzo1MDowNCswMDowMEF1c3RyYWxpYSdzIHNlcnZpY2Ugc2VjdG9yIHJldHVybmVkIHRvIGdyb3d0aCBpbiBGZWJydWFyeSB3aXRoIHRoZSBKdWRvIEJhbmsgQXVzdHJhbGlhIFNlcnZpY2VzICNQTUkgdXAgYXQgNTAuNyAoSmFuIGZpbmFsOiA0OC42KS4gUHJpY2UgcHJlc3N1cmVzIHJlY2VkZWQgaW4gdGhlIHNlY3Rvci4gQEp1ZG9fQmFuayBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9wQ0pGcHkxTHBvIGh0dHBzOi8vdC5jby8zYm0zT3I3V1JJ",
    "Normal|MjAyMy0wMi0yOCAyMzo0NjowNCswMDowMEF1c3RyYWxpYeKAmXMgbWFudWZhY3R1cmluZyBzZWN0b3IgZXhwYW5zaW9uIHJlc3VtZWQgaW4gRmVicnVhcnkgd2l0aCB0aGUgSnVkbyBCYW5rIEF1c3RyYWxpYSBNYW51ZmFjdHVyaW5nICNQTUkgYXQgNTAuNSAoSmFuOiA1MC4wKS4gUHJpY2UgcHJlc3N1cmVzIGVhc2VkIHdpdGggY29zdCBpbmZsYXRpb24gYXQgdGhlIGxvd2VzdCBpbiB0d28geWVhcnMuIEBKdWRvX0JhbmsgUmVhZCBtb3JlOiBodHRwczovL3QuY28vaHJQRHdndzh2UCBodHRwczovL3QuY28vWUhIVEdJZTQyag==",
    "Heading2|LS0tLS0tLS0tLS0tLSBVSyBQTUkgLS0tLS0tLS0tLS0tLS0tLS0tLS0tIDo=",
    "Normal|MjAyMy0wNC0wNSAwODozMzoxMyswMDowMPCfh6zwn4enI1VLIHNlcnZpY2UgcHJvdmlkZXJzIHJlcG9ydGVkIGEgc3VzdGFpbmVkIGV4cGFuc2lvbiBpbiBhY3Rpdml0eSBpbiBNYXJjaCAoI1BNSSBhdCA1Mi45OyBGZWI6IDUzLjUpLiBUaGUgcmVjb3Zlcnkgd2FzIGNlbWVudGVkIGJ5IHRoZSBmYXN0ZXN0IGluY3JlYXNlIGluIG5ldyBvcmRlciB2b2x1bWVzIGluIGEgeWVhciBhbmQgdGhlIHN0cm9uZ2VzdCByaXNlIGluIG5ldyBleHBvcnQgc2FsZXMgaW4gc3VydmV5IGhpc3RvcnkuIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL2xqY1BxejhLSVMgaHR0cHM6Ly90LmNvL0dzNzlqYnRWV1I=",
    "Normal|MjAyMy0wNC0wMyAwODo0MDozMSswMDowMPCfh6zwn4enI1VL4oCZcyBtYW51ZmFjdHVyaW5nIHJlZ2lzdGVyZWQgYSBzdXN0YWluZWQgYW5kIHN0cm9uZ2VyIGNvbnRyYWN0aW9uIGluIE1hcmNoICgjUE1JIGF0IDQ3Ljk7IEZlYjogNDkuMykgYXMgb3V0cHV0IHdhcyBzY2FsZWQgYmFjayBpbiByZXNwb25zZSB0byBzdWJkdWVkIG1hcmtldCBkZW1hbmQuIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL00yZkw2N25sMlMgaHR0cHM6Ly90LmNvL2xEQmR6TnBHZnA=",
    "Normal|MjAyMy0wMy0yNCAwOTozNDowMCswMDowMPCfh6zwn4enQnVveWVkIGJ5IHN0cm9uZyAjc2VydmljZSBzZWN0b3IgcGVyZm9ybWFuY2UsIE1hcmNoIHNhdyBhIHN1c3RhaW5lZCBpbXByb3ZlbWVudCBpbiB0aGUgI1VLIHByaXZhdGUgc2VjdG9yIChoZWFkbGluZSAjUE1JIGF0IDUyLjI7IEZlYjogNTMuMSkgYnV0IG1hbnVmYWN0dXJpbmcgcHJvZHVjdGlvbiBkaXBwZWQsIGhlbGQgYmFjayBieSBzdGlsbCBzdWJkdWVkIG9yZGVyIGJvb2tzLiBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9XcGd0R2p4MDVjIGh0dHBzOi8vdC5jby9tM1JsM3hYMFQ1",
    "Normal|MjAyMy0wMy0yMiAxNTozODoyOCswMDowMPCfh6zwn4enVG9kYXnigJlzICNDUEkgcHJpbnQgZm9yIHRoZSBVSyByb3NlIHVuZXhwZWN0ZWRseSB0byAxMC40JSAoSmFuOiAxMC4xJSksIHJlaWduaXRpbmcgaW5mbGF0aW9uYXJ5IGNvbmNlcm5zIGFuZCBhZGRpbmcgdG8gZXhwZWN0YXRpb25zIHRoYXQgdGhlIEJvRSB3aWxsIHJhaXNlIGludGVyZXN0IHJhdGVzIGFnYWluIGluIFRodXJzZGF54oCZcyBtZWV0aW5nLiBPdXIgVUsgQ29tcG9zaXRlIE91dHB1dCBQcmljZXMgUE1JIHNpZ25hbGxlZCB0aGF0IHByaWNlIHByZXNzdXJlcyByZW1haW4gc3RpY2t5IGh0dHBzOi8vdC5jby9ON2ZrbU9tWGZ4",
    "Normal|MjAyMy0wMy0wOSAwOToyODoxMCswMDowMCNQTUkgI0ZhY3RvZnRoZVdlZWsgVGhlIHJhdGlvIGJldHdlZW4gdGhlIPCfh6zwn4enICNVSyBNYW51ZmFjdHVyaW5nIE5ldyBPcmRlcnMgSW5kZXggYW5kIFF1YW50aXR5IG9mIFB1cmNoYXNlcyBJbmRleCB3YXMgdGhlIGhpZ2hlc3Qgb24gcmVjb3JkIGluIEZlYnJ1YXJ5ICgxLjE3KSwgc2lnbmFsbGluZyBhIHJlbHVjdGFuY2UgYW1vbmcgZmlybXMgdG8gYnVpbGQgaW5wdXQgc3RvY2tzIGluIGxpbmUgd2l0aCBzdGFiaWxpc2luZyBkZW1hbmQuIGh0dHBzOi8vdC5jby9JMXBMMlE5QWRJ",
    "Heading2|LS0tLS0tLS0tLS0tLSBHZXJtYW55IFBNSSAtLS0tLS0tLS0tLS0tLS0tLS0tLS0gOg==",
    "Normal|MjAyMy0wNC0wNSAwODoxMDozMSswMDowMPCfh6nwn4eqU3VwcG9ydGVkIGJ5IGEgc3RyZW5ndGhlbmluZyBpbiB1bmRlcmx5aW5nIGRlbWFuZCwgI0dlcm1hbnnigJlzICNzZXJ2aWNlIHNlY3RvciByZW1haW5lZCBvbiBhbiB1cHdhcmQgdHJhamVjdG9yeSBpbiBNYXJjaCBhbmQgcmVnaXN0ZXJlZCBzb2xpZCBncm93dGggaW4gYWN0aXZpdHkgKCNQTUkgYXQgNTMuNzsgRmViOiA1MC45KS4gUmF0ZXMgb2YgaW5mbGF0aW9uIGNvb2xlZCBidXQgcmVtYWluZWQgc3Vic3RhbnRpYWwuIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL0dFb0VNQVRGWEggaHR0cHM6Ly90LmNvL1MyZlVFd2ozbDc=",
    "Normal|MjAyMy0wNC0wMyAwODoxNTowOSswMDowMPCfh6nwn4eqI0dlcm1hbnkncyBtYW51ZmFjdHVyaW5nIGZlbGwgZGVlcGVyIGludG8gY29udHJhY3Rpb24gdGVycml0b3J5IGluIE1hcmNoICgjUE1JIGF0IDQ0Ljc7IEZlYjogNDYuMyksIGxhcmdlbHkgZHJpdmVuIGJ5IHRoZSBncmVhdGVzdCBpbXByb3ZlbWVudCBpbiBzdXBwbGllciBwZXJmb3JtYW5jZSBzZWVuIGluIHRoZSBzZXJpZXMgaGlzdG9yeSAoaW52ZXJ0ZWQgaW4gdGhlIFBNSSBjYWxjdWxhdGlvbikuIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL2twdW95Snp6NE8gaHR0cHM6Ly90LmNvL0JXNXppMjJWR0M=",
    "Normal|MjAyMy0wMy0wMyAwOToyMTozMyswMDowMPCfh6nwn4eqQmFjay10by1iYWNrIG1hcmdpbmFsIGdyb3d0aCBpbiBGZWJydWFyeSAoaGVhZGxpbmUgI1BNSSBhdCA1MC45OyBKYW46IDUwLjcpIGhhcyBtZWFudCB0aGF0ICNHZXJtYW554oCZcyBzZXJ2aWNlIHNlY3RvciBpcyBub3cgbW9yZSBzdGFibGUuIFRoZSBsYXRlc3QgdXB0dXJuIHdhcyBhY2NvbXBhbmllZCBieSBmcmVzaCBuZXcgb3JkZXIgZ3Jvd3RoIGFuZCBzb21lIGVhc2luZyBpbiBpbmZsYXRpb25hcnkgcHJlc3N1cmVzLiBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9EUHVYR0haekRyIGh0dHBzOi8vdC5jby9aZFc1ZjVCb1k0",
    "Normal|MjAyMy0wMi0yMSAwODo0ODo0NSswMDowMPCfh6nwn4eqI0dlcm1hbnnigJlzIHByaXZhdGUgc2VjdG9yIHJlY29yZGVkIGZyZXNoIGdyb3d0aCBpbiBhY3Rpdml0eSBhbWlkIHNpZ25zIG9mIHJlZHVjZWQgc3RyYWluIG9uIGRlbWFuZCBhbmQgYSByYXBpZCBlYXNpbmcgb2Ygc3VwcGx5LXNpZGUgY29uc3RyYWludHMuIFRoZSBoZWFkbGluZSAjUE1JIHJlYWRpbmcgcm9zZSB0byA1MS4xIGluIEZlYnJ1YXJ5IChKYW46IDQ5LjkpIHRvIHJlYWNoIGFuIDgtbW9udGggaGlnaC4gUmVhZCBtb3JlOiBodHRwczovL3QuY28vVjdlQ1Q5a1NhVSBodHRwczovL3QuY28vanNheUZsZnk3ZA==",
    "Normal|MjAyMy0wMi0wMyAwOTozNjozNSswMDowMPCfh6nwn4eqI0dlcm1hbnnigJlzIHNlcnZpY2Ugc2VjdG9yIHNpZ25hbGxlZCBhIHJlbmV3ZWQgdXB0dXJuIGluIGJ1c2luZXNzIGFjdGl2aXR5IGluIHRoZSBmaXJzdCBtb250aCBvZiAyMDIzIHdpdGggdGhlICNQTUkgcmVhZGluZyByaXNpbmcgdG8gNTAuNyAoRGVjOiA0OS4yKSBidXQgZGVtYW5kIHJlbWFpbmVkIGZyYWdpbGUgYXMgbmV3IGJ1c2luZXNzIGZlbGwgYWdhaW4uIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL2oyM0tKM3g1YW4gaHR0cHM6Ly90LmNvL0tlQkxPcHFabEQ=",
    "Heading2|LS0tLS0tLS0tLS0tLSBKYXBhbiBQTUkgLS0tLS0tLS0tLS0tLS0tLS0tLS0tIDo=",
    "Normal|MjAyMy0wNC0wNSAwMTozOTowNCswMDowMEphcGFu4oCZcyBzZXJ2aWNlIHNlY3RvciBleHBhbnNpb24gYWNjZWxlcmF0ZWQgaW4gTWFyY2ggd2l0aCB0aGUgYXUgSmlidW4gQmFuayBKYXBhbiBTZXJ2aWNlcyAjUE1JIGF0IDU1LjAgKEZlYjogNTQuMCkuIE5ldyBidXNpbmVzcyByb3NlIGF0IHNoYXJwZXN0IHBhY2Ugc2luY2UgRmViIDE5LiBAamlidW5iYW5rIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL2ZSaldCcm1kTHEgaHR0cHM6Ly90LmNvLzc5S2VucXBBcjE=",
    "Normal|MjAyMy0wNC0wMyAwMjowOToxNCswMDowMEphcGFuJ3MgbWFudWZhY3R1cmluZyBzZWN0b3Igc2F3IGEgc29mdGVyIGNvbnRyYWN0aW9uIGluIE1hcmNoIHdpdGggdGhlIGF1IEppYnVuIEJhbmsgSmFwYW4gTWFudWZhY3R1cmluZyAjUE1JIHVwIGF0IDQ5LjIgKEZlYjogNDcuNykuIElucHV0IHByaWNlIGluZmxhdGlvbiBlYXNlZCB0byBzb2Z0ZXN0IHNpbmNlIEF1Z3VzdCAyMDIxLiBAamlidW5iYW5rIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL2hWd3NHNnRQdU8gaHR0cHM6Ly90LmNvL1V0WmxZRnA2Q3U=",
    "Normal|MjAyMy0wMy0yNCAwMDozNDowNyswMDowMEphcGFu4oCZcyBwcml2YXRlIHNlY3RvciBvdXRwdXQgaW5jcmVhc2VkIGF0IHRoZSBmYXN0ZXN0IHJhdGUgaW4gbmluZSBtb250aHMgd2l0aCB0aGUgYXUgSmlidW4gQmFuayBGbGFzaCBKYXBhbiBDb21wb3NpdGUgI1BNSSBwb3N0aW5nIDUxLjkgKEZlYiBmaW5hbDogNTEuMSkuIFNlcnZpY2VzIGFjdGl2aXR5IGdyb3d0aCBhY2NlbGVyYXRlZC4gQGppYnVuYmFuayBSZWFkIG1vcmU6IGh0dHBzOi8vdC5jby9ENVhRN1QyZUlOIGh0dHBzOi8vdC5jby80YlZwakg0VlY1",
    "Normal|MjAyMy0wMy0wMyAwMDo0NzozOSswMDowMEphcGFuJ3Mgc2VydmljZSBzZWN0b3IgYWN0aXZpdHkgZXhwYW5kZWQgc29saWRseSBpbiBGZWJydWFyeSB3aXRoIGF1IEppYnVuIEJhbmsgSmFwYW4gU2VydmljZXMgI1BNSSBhdCA1NC4wIChKYW4gZmluYWw6IDUyLjMpLiBCdXNpbmVzcyBjb25maWRlbmNlIHJvc2UgdG8gYSBmb3VyLW1vbnRoIGhpZ2guIEBqaWJ1bmJhbmsgUmVhZCBtb3JlOiBodHRwczovL3QuY28veUR1RGtBNjh3TCBodHRwczovL3QuY28vS2xaUVQyUnloZg==",
    "Normal|MjAyMy0wMy0wMSAwMDozNToyMiswMDowMEphcGFu4oCZcyBtYW51ZmFjdHVyaW5nIHNlY3RvciBzYXcgZnVydGhlciBmYWxscyBpbiBvdXRwdXQgYW5kIG5ldyBvcmRlcnMgaW4gRmVicnVhcnkuIFRoZSBhdSBKaWJ1biBCYW5rIEphcGFuIE1hbnVmYWN0dXJpbmcgI1BNSSBlYXNlZCB0byA0Ny43IChKYW4gZmluYWw6IDQ4LjkpLiBAamlidW5iYW5rIFJlYWQgbW9yZTogaHR0cHM6Ly90LmNvL0pHOG95Y1h2UnIgaHR0cHM6Ly90LmNvL0YwaUpzMDhLMkI="
)

# Remove every paragraph from the first section heading (paragraph 3,
# "China PMI") through the end of the body, keeping the Title and the
# "PMI" Heading1 paragraph intact.
$startPara = $d.Paragraphs.Item(3)
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# Re-append the new paragraphs, one at a time, after the last remaining
# paragraph ("PMI" Heading1).
foreach ($item in $items) {
    $sep = $item.IndexOf("|")
    $styleToken = $item.Substring(0, $sep)
    $b64 = $item.Substring($sep + 1)
    $bytes = [System.Convert]::FromBase64String($b64)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)

    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

    if ($styleToken -eq "Heading2") {
        $newPara.Style = "Heading 2"
    } else {
        $newPara.Style = "Normal"
    }

    $r = $newPara.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
